# Updates the "want to go" counts (column F) and a couple of sold-out status
# labels (column G) across the four worksheets, and appends a newly scraped
# local-life event as row 4 of the "本地生活" sheet.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 898
$ws1.Cells.Item(5, 6).Value = 828
$ws1.Cells.Item(6, 6).Value = 0
$ws1.Cells.Item(7, 6).Value = 0
$ws1.Cells.Item(8, 6).Value = 1476
$ws1.Cells.Item(9, 6).Value = 38297
$ws1.Cells.Item(9, 7).Value = "已售罄"
$ws1.Cells.Item(10, 6).Value = 7960
$ws1.Cells.Item(11, 6).Value = 0
$ws1.Cells.Item(13, 6).Value = 651
$ws1.Cells.Item(15, 6).Value = 58
$ws1.Cells.Item(17, 6).Value = 0
$ws1.Cells.Item(18, 6).Value = 0
$ws1.Cells.Item(20, 6).Value = 0
$ws1.Cells.Item(22, 6).Value = 185
$ws1.Cells.Item(23, 6).Value = 0
$ws1.Cells.Item(25, 6).Value = 500
$ws1.Cells.Item(26, 6).Value = 359
$ws1.Cells.Item(27, 6).Value = 0
$ws1.Cells.Item(28, 6).Value = 532
$ws1.Cells.Item(29, 6).Value = 31
$ws1.Cells.Item(30, 6).Value = 338
$ws1.Cells.Item(33, 6).Value = 345
$ws1.Cells.Item(34, 6).Value = 164
$ws1.Cells.Item(35, 6).Value = 189
$ws1.Cells.Item(36, 6).Value = 0
$ws1.Cells.Item(37, 6).Value = 163
$ws1.Cells.Item(38, 6).Value = 41
$ws1.Cells.Item(40, 6).Value = 316

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 3
$ws2.Cells.Item(3, 6).Value = 0
$ws2.Cells.Item(4, 6).Value = 321
$ws2.Cells.Item(5, 6).Value = 4350
$ws2.Cells.Item(7, 6).Value = 267
$ws2.Cells.Item(8, 6).Value = 0
$ws2.Cells.Item(10, 6).Value = 0
$ws2.Cells.Item(11, 6).Value = 78
$ws2.Cells.Item(14, 6).Value = 0
$ws2.Cells.Item(15, 6).Value = 156
$ws2.Cells.Item(16, 6).Value = 0
$ws2.Cells.Item(17, 6).Value = 4350

# --- 本地生活 (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 0

# Append a new row (row 4) for a newly scraped event, copying the formatting
# of the row above it first so styles (bold index column, borders, etc.)
# carry over correctly.
$ws3.Range("A3:I3").Copy()
$ws3.Range("A4:I4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Cells.Item(4, 1).Value = 3
# Column B holds plain text dates (e.g. "2024-07-12"); force text formatting
# first so Excel does not auto-convert the string into a date serial number,
# then restore the original (unformatted) style copied from the row above.
$ws3.Cells.Item(4, 2).NumberFormat = "@"
$ws3.Cells.Item(4, 2).Value = "2024-07-12"
$ws3.Cells.Item(3, 2).Copy()
$ws3.Cells.Item(4, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws3.Cells.Item(4, 3).Value = "广州·全职高手×HAPPY ZOO 全职高手 十周年咖啡厅"
$ws3.Cells.Item(4, 4).Value = "多宝街道恩宁路十一甫新街7号 啡约咖啡馆"
$ws3.Cells.Item(4, 5).Value = "2024.07.12 00:00-08.18 23:59"
$ws3.Cells.Item(4, 6).Value = 0
$ws3.Cells.Item(4, 7).Value = 10
$ws3.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=88806"
$ws3.Cells.Item(4, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/gNQIvlhI1720418693552.png"

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 1667
$ws4.Cells.Item(4, 6).Value = 266
$ws4.Cells.Item(7, 6).Value = 828
$ws4.Cells.Item(8, 6).Value = 0
$ws4.Cells.Item(9, 6).Value = 1476
$ws4.Cells.Item(10, 6).Value = 38297
$ws4.Cells.Item(10, 7).Value = "已售罄"
$ws4.Cells.Item(13, 6).Value = 5
$ws4.Cells.Item(16, 6).Value = 0
$ws4.Cells.Item(18, 6).Value = 0
$ws4.Cells.Item(24, 6).Value = 78
$ws4.Cells.Item(25, 6).Value = 168
$ws4.Cells.Item(26, 6).Value = 0
$ws4.Cells.Item(27, 6).Value = 6
$ws4.Cells.Item(30, 6).Value = 0
$ws4.Cells.Item(31, 6).Value = 185
$ws4.Cells.Item(32, 6).Value = 952
$ws4.Cells.Item(33, 6).Value = 0
$ws4.Cells.Item(34, 6).Value = 0
$ws4.Cells.Item(37, 6).Value = 532
$ws4.Cells.Item(40, 6).Value = 0
$ws4.Cells.Item(42, 6).Value = 65
$ws4.Cells.Item(44, 6).Value = 0
$ws4.Cells.Item(45, 6).Value = 0
$ws4.Cells.Item(46, 6).Value = 941
$ws4.Cells.Item(47, 6).Value = 316
$ws4.Cells.Item(49, 6).Value = 13
$ws4.Cells.Item(50, 6).Value = 44

Write-Output "Edit applied successfully"
